# SSU - Brisanje Predviđanja: apply the commit's content edits via Word COM
# automation (Find/Replace + bookmark-targeted range edits).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: "... funkcionalnosti brisanje predviđanja" -> "... brisanjepredviđanja"
#    (remove the space run between "brisanje" and "predviđanja")
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$null = $rng.Find.Execute(
    "funkcionalnosti brisanje predviđanja",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "funkcionalnosti brisanjepredviđanja", 1)

# ---------------------------------------------------------------------------
# 2) "... projektnog tima u razvoju projekta i testiranju ..." ->
#    "... projektnog tima u razvojuprojekta i testiranju ..."
#    (remove the space run between "ima u razvoju" and "projekta i testiranju")
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$null = $rng.Find.Execute(
    "tima u razvoju projekta i testiranju",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "tima u razvojuprojekta i testiranju", 1)

# ---------------------------------------------------------------------------
# 3) Heading "Scenario brisanje ideja" -> "Scenario brisanjeideja"
#    (remove the space run between "brisanje" and "ideja"); use the
#    heading's own bookmark so the cached TOC entry is left untouched.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_Toc34567528")
$hdgRng = $bm.Range
$null = $hdgRng.Find.Execute(
    "Scenario brisanje ideja",
    $false, $false, $false, $false, $false,
    $true, 0, $false,
    "Scenario brisanjeideja", 1)

# ---------------------------------------------------------------------------
# 4) "... bira taster sa natpisom “Odustani” koji ga vraća ..." ->
#    "... bira taster X koji ga vraća ..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$null = $rng.Find.Execute(
    "bira taster sa natpisom " + [char]0x201C + "Odustani" + [char]0x201D + " koji",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "bira taster X koji", 1)
